$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.081981583213166
$ws.Cells.Item(2, 4).Value = 1.082531993875297
$ws.Cells.Item(2, 5).Value = 1.084519647350569
$ws.Cells.Item(2, 6).Value = 1.094501719293237
$ws.Cells.Item(2, 9).Value = 1.057984687742426
$ws.Cells.Item(2, 10).Value = 1.086853476581865
$ws.Cells.Item(2, 11).Value = 1.085199966538849
$ws.Cells.Item(2, 12).Value = 1.087182456347876
$ws.Cells.Item(2, 13).Value = 1.097138901975255
$ws.Cells.Item(2, 14).Value = 1.088396933182447
# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.083698833586853
$ws.Cells.Item(3, 4).Value = 1.083902060528331
$ws.Cells.Item(3, 5).Value = 1.086032015474934
$ws.Cells.Item(3, 6).Value = 1.09600629205567
$ws.Cells.Item(3, 9).Value = 1.058499236461373
$ws.Cells.Item(3, 10).Value = 1.088229113723273
$ws.Cells.Item(3, 11).Value = 1.086387350657828
$ws.Cells.Item(3, 12).Value = 1.088512176327696
$ws.Cells.Item(3, 13).Value = 1.098462717987029
$ws.Cells.Item(3, 14).Value = 1.089774523886384
# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.084807553590431
$ws.Cells.Item(4, 4).Value = 1.084786300694272
$ws.Cells.Item(4, 5).Value = 1.087008123024208
$ws.Cells.Item(4, 6).Value = 1.096977649583864
$ws.Cells.Item(4, 9).Value = 1.058829573980172
$ws.Cells.Item(4, 10).Value = 1.089116411017484
$ws.Cells.Item(4, 11).Value = 1.087152843556756
$ws.Cells.Item(4, 12).Value = 1.089369587268702
$ws.Cells.Item(4, 13).Value = 1.099316596634084
$ws.Cells.Item(4, 14).Value = 1.090663081244527
# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.085273084585707
$ws.Cells.Item(5, 4).Value = 1.085157497810774
$ws.Cells.Item(5, 5).Value = 1.087417892264753
$ws.Cells.Item(5, 6).Value = 1.097385491655685
$ws.Cells.Item(5, 9).Value = 1.058967827509609
$ws.Cells.Item(5, 10).Value = 1.089488763336931
$ws.Cells.Item(5, 11).Value = 1.08747398963343
$ws.Cells.Item(5, 12).Value = 1.089729333389368
$ws.Cells.Item(5, 13).Value = 1.099674926118227
$ws.Cells.Item(5, 14).Value = 1.091035962346978
# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.085351215953054
$ws.Cells.Item(6, 4).Value = 1.08521979219776
$ws.Cells.Item(6, 5).Value = 1.08748666025954
$ws.Cells.Item(6, 6).Value = 1.097453940140581
$ws.Cells.Item(6, 9).Value = 1.058991004635804
$ws.Cells.Item(6, 10).Value = 1.08955124409314
$ws.Cells.Item(6, 11).Value = 1.087527872609877
$ws.Cells.Item(6, 12).Value = 1.089789695044728
$ws.Cells.Item(6, 13).Value = 1.099735053904138
$ws.Cells.Item(6, 14).Value = 1.091098531833026
# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.08481377627938
$ws.Cells.Item(7, 4).Value = 1.084791262746999
$ws.Cells.Item(7, 5).Value = 1.087013600669566
$ws.Cells.Item(7, 6).Value = 1.096983101205163
$ws.Cells.Item(7, 9).Value = 1.058831423762022
$ws.Cells.Item(7, 10).Value = 1.089121389015583
$ws.Cells.Item(7, 11).Value = 1.087157137335391
$ws.Cells.Item(7, 12).Value = 1.089374396987245
$ws.Cells.Item(7, 13).Value = 1.099321387155737
$ws.Cells.Item(7, 14).Value = 1.090668066311954
# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.082562452383676
$ws.Cells.Item(8, 4).Value = 1.082995493081082
$ws.Cells.Item(8, 5).Value = 1.085031282717107
$ws.Cells.Item(8, 6).Value = 1.09501065888311
$ws.Cells.Item(8, 9).Value = 1.058159124982624
$ws.Cells.Item(8, 10).Value = 1.087318972309716
$ws.Cells.Item(8, 11).Value = 1.085601838695971
$ws.Cells.Item(8, 12).Value = 1.087632470299056
$ws.Cells.Item(8, 13).Value = 1.097586860665499
$ws.Cells.Item(8, 14).Value = 1.088863089967625
# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.078575936704114
$ws.Cells.Item(9, 4).Value = 1.079813193107303
$ws.Cells.Item(9, 5).Value = 1.081518570230104
$ws.Cells.Item(9, 6).Value = 1.09151763081623
$ws.Cells.Item(9, 9).Value = 1.056954258919808
$ws.Cells.Item(9, 10).Value = 1.084120721523485
$ws.Cells.Item(9, 11).Value = 1.082839172301233
$ws.Cells.Item(9, 12).Value = 1.084539476343338
$ws.Cells.Item(9, 13).Value = 1.094509139106293
$ws.Cells.Item(9, 14).Value = 1.085660297298433
# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.075904395831092
$ws.Cells.Item(10, 4).Value = 1.07767900767124
$ws.Cells.Item(10, 5).Value = 1.079162871164242
$ws.Cells.Item(10, 6).Value = 1.089176630618242
$ws.Cells.Item(10, 9).Value = 1.056137157497866
$ws.Cells.Item(10, 10).Value = 1.081972989106421
$ws.Cells.Item(10, 11).Value = 1.08098201565619
$ws.Cells.Item(10, 12).Value = 1.082461025342585
$ws.Cells.Item(10, 13).Value = 1.092442407767931
$ws.Cells.Item(10, 14).Value = 1.083509514854989
# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.074744116751408
$ws.Cells.Item(11, 4).Value = 1.076751746259069
$ws.Cells.Item(11, 5).Value = 1.078139375547861
$ws.Cells.Item(11, 6).Value = 1.088159882689998
$ws.Cells.Item(11, 9).Value = 1.055779993214888
$ws.Cells.Item(11, 10).Value = 1.081039156644373
$ws.Cells.Item(11, 11).Value = 1.080174071710369
$ws.Cells.Item(11, 12).Value = 1.081556983143011
$ws.Cells.Item(11, 13).Value = 1.091543813306838
$ws.Cells.Item(11, 14).Value = 1.082574356243732
# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.074312598025476
$ws.Cells.Item(12, 4).Value = 1.076406835880647
$ws.Cells.Item(12, 5).Value = 1.077758669888696
$ws.Cells.Item(12, 6).Value = 1.087781741742255
$ws.Cells.Item(12, 9).Value = 1.055646817048754
$ws.Cells.Item(12, 10).Value = 1.080691698337805
$ws.Cells.Item(12, 11).Value = 1.079873385991208
$ws.Cells.Item(12, 12).Value = 1.081220559007932
$ws.Cells.Item(12, 13).Value = 1.091209469450847
$ws.Cells.Item(12, 14).Value = 1.082226404506515
# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.074405184924378
$ws.Cells.Item(13, 4).Value = 1.076480842428022
$ws.Cells.Item(13, 5).Value = 1.077840356921591
$ws.Cells.Item(13, 6).Value = 1.087862875984174
$ws.Cells.Item(13, 9).Value = 1.055675406931878
$ws.Cells.Item(13, 10).Value = 1.080766256292616
$ws.Cells.Item(13, 11).Value = 1.079937910494474
$ws.Cells.Item(13, 12).Value = 1.081292751507579
$ws.Cells.Item(13, 13).Value = 1.091281213144393
$ws.Cells.Item(13, 14).Value = 1.082301068342174
# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.074708458375615
$ws.Cells.Item(14, 4).Value = 1.076723245805315
$ws.Cells.Item(14, 5).Value = 1.078107917268357
$ws.Cells.Item(14, 6).Value = 1.088128635216106
$ws.Cells.Item(14, 9).Value = 1.055768995256385
$ws.Cells.Item(14, 10).Value = 1.081010447770903
$ws.Cells.Item(14, 11).Value = 1.080149228818107
$ws.Cells.Item(14, 12).Value = 1.081529186994383
$ws.Cells.Item(14, 13).Value = 1.091516187953057
$ws.Cells.Item(14, 14).Value = 1.08254560660037
# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.074895243200795
$ws.Cells.Item(15, 4).Value = 1.07687253401262
$ws.Cells.Item(15, 5).Value = 1.078272698882574
$ws.Cells.Item(15, 6).Value = 1.088292314855377
$ws.Cells.Item(15, 9).Value = 1.055826590437885
$ws.Cells.Item(15, 10).Value = 1.081160823451746
$ws.Cells.Item(15, 11).Value = 1.080279351908304
$ws.Cells.Item(15, 12).Value = 1.081674779815378
$ws.Cells.Item(15, 13).Value = 1.091660888342501
$ws.Cells.Item(15, 14).Value = 1.082696195831927
# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.075981324965856
$ws.Cells.Item(16, 4).Value = 1.07774047969942
$ws.Cells.Item(16, 5).Value = 1.0792307230887
$ws.Cells.Item(16, 6).Value = 1.089244042856405
$ws.Cells.Item(16, 9).Value = 1.056160790178145
$ws.Cells.Item(16, 10).Value = 1.082034882274292
$ws.Cells.Item(16, 11).Value = 1.081035555595466
$ws.Cells.Item(16, 12).Value = 1.082520937037038
$ws.Cells.Item(16, 13).Value = 1.092501965833303
$ws.Cells.Item(16, 14).Value = 1.083571495918257
# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.07666165132398
$ws.Cells.Item(17, 4).Value = 1.078284068863473
$ws.Cells.Item(17, 5).Value = 1.079830729804648
$ws.Cells.Item(17, 6).Value = 1.089840203214262
$ws.Cells.Item(17, 9).Value = 1.056369522990963
$ws.Cells.Item(17, 10).Value = 1.082582116414154
$ws.Cells.Item(17, 11).Value = 1.081508881658501
$ws.Cells.Item(17, 12).Value = 1.083050613375337
$ws.Cells.Item(17, 13).Value = 1.093028556159377
$ws.Cells.Item(17, 14).Value = 1.084119507193366
# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.077058139083923
$ws.Cells.Item(18, 4).Value = 1.078600832572594
$ws.Cells.Item(18, 5).Value = 1.080180370429751
$ws.Cells.Item(18, 6).Value = 1.090187637223066
$ws.Cells.Item(18, 9).Value = 1.056490950077549
$ws.Cells.Item(18, 10).Value = 1.082900938224487
$ws.Cells.Item(18, 11).Value = 1.08178460050293
$ws.Cells.Item(18, 12).Value = 1.083359173836668
$ws.Cells.Item(18, 13).Value = 1.093335352548287
$ws.Cells.Item(18, 14).Value = 1.084438781767239
# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.077193274805965
$ws.Cells.Item(19, 4).Value = 1.078708789828107
$ws.Cells.Item(19, 5).Value = 1.080299532751196
$ws.Cells.Item(19, 6).Value = 1.090306053380536
$ws.Cells.Item(19, 9).Value = 1.056532298917381
$ws.Cells.Item(19, 10).Value = 1.08300958575015
$ws.Cells.Item(19, 11).Value = 1.081878552054102
$ws.Cells.Item(19, 12).Value = 1.083464319080132
$ws.Cells.Item(19, 13).Value = 1.093439902364587
$ws.Cells.Item(19, 14).Value = 1.084547583584851
# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.076588693519884
$ws.Cells.Item(20, 4).Value = 1.078225778280041
$ws.Cells.Item(20, 5).Value = 1.079766389319281
$ws.Cells.Item(20, 6).Value = 1.089776271584409
$ws.Cells.Item(20, 9).Value = 1.056347161409056
$ws.Cells.Item(20, 10).Value = 1.082523441796599
$ws.Cells.Item(20, 11).Value = 1.081458136006217
$ws.Cells.Item(20, 12).Value = 1.082993824630011
$ws.Cells.Item(20, 13).Value = 1.092972094790804
$ws.Cells.Item(20, 14).Value = 1.084060749251125
# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.074619166922451
$ws.Cells.Item(21, 4).Value = 1.076651877513402
$ws.Cells.Item(21, 5).Value = 1.078029142261114
$ws.Cells.Item(21, 6).Value = 1.088050388976091
$ws.Cells.Item(21, 9).Value = 1.055741449935612
$ws.Cells.Item(21, 10).Value = 1.080938555895387
$ws.Cells.Item(21, 11).Value = 1.080087016921226
$ws.Cells.Item(21, 12).Value = 1.081459579926716
$ws.Cells.Item(21, 13).Value = 1.091447009446575
$ws.Cells.Item(21, 14).Value = 1.082473612630144
# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.07337771847966
$ws.Cells.Item(22, 4).Value = 1.075659496295581
$ws.Cells.Item(22, 5).Value = 1.07693377073576
$ws.Cells.Item(22, 6).Value = 1.086962500540485
$ws.Cells.Item(22, 9).Value = 1.055357664784211
$ws.Cells.Item(22, 10).Value = 1.079938647357865
$ws.Cells.Item(22, 11).Value = 1.07922158236214
$ws.Cells.Item(22, 12).Value = 1.080491330745653
$ws.Cells.Item(22, 13).Value = 1.090484847784768
$ws.Cells.Item(22, 14).Value = 1.081472284107807
# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.074036135736874
$ws.Cells.Item(23, 4).Value = 1.076185846396259
$ws.Cells.Item(23, 5).Value = 1.077514745896876
$ws.Cells.Item(23, 6).Value = 1.087539476491917
$ws.Cells.Item(23, 9).Value = 1.055561398078669
$ws.Cells.Item(23, 10).Value = 1.080469047058945
$ws.Cells.Item(23, 11).Value = 1.079680687568487
$ws.Cells.Item(23, 12).Value = 1.081004964307421
$ws.Cells.Item(23, 13).Value = 1.090995222848946
$ws.Cells.Item(23, 14).Value = 1.082003437037301
# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.076621661046423
$ws.Cells.Item(24, 4).Value = 1.078252118220781
$ws.Cells.Item(24, 5).Value = 1.079795463044438
$ws.Cells.Item(24, 6).Value = 1.089805160451384
$ws.Cells.Item(24, 9).Value = 1.056357266644129
$ws.Cells.Item(24, 10).Value = 1.082549955475729
$ws.Cells.Item(24, 11).Value = 1.081481066907608
$ws.Cells.Item(24, 12).Value = 1.083019486227006
$ws.Cells.Item(24, 13).Value = 1.092997608349302
$ws.Cells.Item(24, 14).Value = 1.08408730058272
# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.079608929305719
$ws.Cells.Item(25, 4).Value = 1.080638079538764
$ws.Cells.Item(25, 5).Value = 1.082429087022794
$ws.Cells.Item(25, 6).Value = 1.092422783994747
$ws.Cells.Item(25, 9).Value = 1.057268167207224
$ws.Cells.Item(25, 10).Value = 1.084950240377431
$ws.Cells.Item(25, 11).Value = 1.083556056581037
$ws.Cells.Item(25, 12).Value = 1.085341941167817
$ws.Cells.Item(25, 13).Value = 1.095307387289264
$ws.Cells.Item(25, 14).Value = 1.0864909941643

Write-Host "Updated vm_pu values for case with 380 kV"
